# Reorganize Ar4 and Ar3: the "right_arm" table (J1:O5) becomes the "head"
# table (J1:O7) with new columns (Name, Pin, Min, Dir, Max, Dir2) and new
# rows describing head pivot points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table first (while it still has its old name/columns), then
# rename it - doing it in this order keeps the resize effective.
$lo = $ws.ListObjects("right_arm")
$lo.Resize($ws.Range("J1:O7"))
$lo.Name = "head"

# Rewrite the header row for the resized table.
$ws.Range("J1").Value = "Name"
$ws.Range("K1").Value = "Pin"
$ws.Range("L1").Value = "Min"
$ws.Range("M1").Value = "Dir"
$ws.Range("N1").Value = "Max"
$ws.Range("O1").Value = "Dir2"

# Fill in the new body rows (Name / Pin columns only - Min/Dir/Max/Dir2 are
# left blank for these rows). Entered in this particular order so that new
# shared-string entries land in the same order as the authored workbook.
$ws.Range("J3").Value = "side_eye"
$ws.Range("J5").Value = "up_eye"
$ws.Range("J2").Value = "head "
$ws.Range("J6").Value = "side_neck"
$ws.Range("J4").Value = "mouth"
$ws.Range("J7").Value = "neck"

$ws.Range("K2").Value = 8
$ws.Range("K3").Value = 9
$ws.Range("K4").Value = 10
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 13

# Clear the leftover values from the old "right_arm" rows (L2:L5 used to
# hold numeric data for the 4-row table; the new 6-row table leaves those
# columns empty for rows 2-5).
$ws.Range("L2").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("L5").ClearContents()

# Match the saved selection.
$ws.Range("J11").Select() | Out-Null
